$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "BUG LIST"
$ws.Name = "BUG LIST"

# Clear out the old content (rows 1-7, columns A-E) so stale cells (e.g. old row 7) go away
$ws.Range("A1:E7").ClearContents()

# Header row (order chosen so shared strings are appended in the same order
# as the target workbook: DEVICE, BUG , VERSION, STATUS, TEST)
$ws.Range("C1").Value = "DEVICE"
$ws.Range("B1").Value = "BUG "
$ws.Range("D1").Value = "VERSION"
$ws.Range("E1").Value = "STATUS"
$ws.Range("A1").Value = "TEST"

# Data rows
$ws.Range("A2").Value = "CS current"
$ws.Range("B2").Value = "Test Mode entry"
$ws.Range("C2").Value = 4832
$ws.Range("D2").Value = "A0"
$ws.Range("E2").Value = "NOT VERIFIED"

$ws.Range("A3").Value = "CS current"
$ws.Range("B3").Value = "Test Mode entry"
$ws.Range("C3").Value = 4832
$ws.Range("D3").Value = "A1"
$ws.Range("E3").Value = "NOT VERIFIED"

$ws.Range("A4").Value = "CS current"
$ws.Range("B4").Value = "Test Mode entry"
$ws.Range("C4").Value = 4832
$ws.Range("D4").Value = "A2"
$ws.Range("E4").Value = "NOT VERIFIED"

$ws.Range("A5").Value = "CS current"
$ws.Range("B5").Value = "Test Mode entry"
$ws.Range("C5").Value = 4832
$ws.Range("D5").Value = "A3"
$ws.Range("E5").Value = "NOT VERIFIED"

# Move the active selection, matching the target workbook's saved cursor position
$ws.Range("I13").Select()
